# Auto-generated edit script applying the recorded Sheets diff
# to the Brynhildr_Profits workbook (profession sheets ALC/ARM/BSM/CRP/CUL/GSM/LTW).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 29
$ws.Range("H29").Value = 3371.4285
$ws.Range("I29").Value = 720
$ws.Range("K29").Value = 2160
$ws.Range("M29").Value = -1879
# Row 70
$ws.Range("H70").Value = 2719.3333
$ws.Range("J70").Value = 2719.3333
$ws.Range("L70").Value = 8157.999899999999
$ws.Range("N70").Value = -8697.999899999999
# Row 73
$ws.Range("H73").Value = 2719.3333
$ws.Range("J73").Value = 2719.3333
$ws.Range("L73").Value = 8157.999899999999
$ws.Range("N73").Value = -10029.9999
# Row 86
$ws.Range("H86").Value = 5807.8184
$ws.Range("I86").Value = 7666
$ws.Range("J86").Value = 5111
$ws.Range("K86").Value = 7666
$ws.Range("L86").Value = 5111
$ws.Range("M86").Value = -6543
$ws.Range("N86").Value = -7357
# Row 89
$ws.Range("H89").Value = 5807.8184
$ws.Range("I89").Value = 7666
$ws.Range("J89").Value = 5111
$ws.Range("K89").Value = 38330
$ws.Range("L89").Value = 25555
$ws.Range("M89").Value = -32714
$ws.Range("N89").Value = -36787
# Row 137
$ws.Range("H137").Value = 25647250
$ws.Range("I137").Value = 100002240
$ws.Range("J137").Value = 7599.793
$ws.Range("K137").Value = 300006720
$ws.Range("L137").Value = 22799.379
$ws.Range("M137").Value = -300004170
$ws.Range("N137").Value = -27899.379
# Row 138
$ws.Range("H138").Value = 2966.25
$ws.Range("I138").Value = 2136.2307
$ws.Range("J138").Value = 3314.3225
$ws.Range("K138").Value = 6408.6921
$ws.Range("L138").Value = 9942.967500000001
$ws.Range("M138").Value = -1268.6921
$ws.Range("N138").Value = -20222.9675

$ws = $wb.Worksheets.Item("ARM")
# Row 17
$ws.Range("H17").Value = 99
$ws.Range("J17").Value = 99
$ws.Range("L17").Value = 99
$ws.Range("N17").Value = -445
# Row 32
$ws.Range("H32").Value = 186449.11
$ws.Range("I32").Value = 299429.38
$ws.Range("J32").Value = 19434.783
$ws.Range("K32").Value = 299429.38
$ws.Range("L32").Value = 19434.783
$ws.Range("M32").Value = -299142.38
$ws.Range("N32").Value = -20008.783
# Row 45
$ws.Range("H45").Value = 1853.8889
$ws.Range("I45").Value = 1671.25
$ws.Range("K45").Value = 1671.25
$ws.Range("M45").Value = -1294.25
# Row 94
$ws.Range("H94").Value = 50329.5
$ws.Range("J94").Value = 50329.5
$ws.Range("L94").Value = 50329.5
$ws.Range("N94").Value = -52131.5
# Row 122
$ws.Range("H122").Value = 1864.1875
$ws.Range("I122").Value = 1612.3334
$ws.Range("K122").Value = 4837.0002
$ws.Range("M122").Value = -2387.0002

$ws = $wb.Worksheets.Item("BSM")
# Row 46
$ws.Range("H46").Value = 4999.3335
$ws.Range("J46").Value = 4999.3335
$ws.Range("L46").Value = 4999.3335
$ws.Range("N46").Value = -5595.3335
# Row 134
$ws.Range("H134").Value = 3479450
$ws.Range("I134").Value = 5117.0938
$ws.Range("K134").Value = 15351.2814
$ws.Range("M134").Value = -12816.2814

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3476127.5
$ws.Range("I31").Value = 5054366.5
$ws.Range("J31").Value = 4002
$ws.Range("K31").Value = 5054366.5
$ws.Range("L31").Value = 4002
$ws.Range("M31").Value = -5054071.5
$ws.Range("N31").Value = -4592
# Row 34
$ws.Range("H34").Value = 3476127.5
$ws.Range("I34").Value = 5054366.5
$ws.Range("J34").Value = 4002
$ws.Range("K34").Value = 5054366.5
$ws.Range("L34").Value = 4002
$ws.Range("M34").Value = -5054164.5
$ws.Range("N34").Value = -4406
# Row 58
$ws.Range("H58").Value = 3631841
$ws.Range("J58").Value = 10434003
$ws.Range("L58").Value = 10434003
$ws.Range("N58").Value = -10434409
# Row 134
$ws.Range("H134").Value = 1944.6976
$ws.Range("I134").Value = 1576.3243
$ws.Range("J134").Value = 4216.3335
$ws.Range("K134").Value = 4728.9729
$ws.Range("L134").Value = 12649.0005
$ws.Range("M134").Value = -2193.9729
$ws.Range("N134").Value = -17719.0005
# Row 136
$ws.Range("H136").Value = 3631841
$ws.Range("J136").Value = 10434003
$ws.Range("L136").Value = 31302009
$ws.Range("N136").Value = -31307109

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 8124.8887
$ws.Range("I3").Value = 2624.8
$ws.Range("K3").Value = 7874.400000000001
$ws.Range("M3").Value = -7762.400000000001
# Row 18
$ws.Range("H18").Value = 84412
$ws.Range("I18").Value = 125931.75
$ws.Range("J18").Value = 1372.5
$ws.Range("K18").Value = 377795.25
$ws.Range("L18").Value = 4117.5
$ws.Range("M18").Value = -377626.25
$ws.Range("N18").Value = -4455.5
# Row 60
$ws.Range("H60").Value = 539.6923
$ws.Range("I60").Value = 450.18182
$ws.Range("J60").Value = 1032
$ws.Range("K60").Value = 1350.54546
$ws.Range("L60").Value = 3096
$ws.Range("M60").Value = -1099.54546
$ws.Range("N60").Value = -3598
# Row 62
$ws.Range("H62").Value = 2824.75
$ws.Range("I62").Value = 2766.3333
$ws.Range("K62").Value = 8298.999899999999
$ws.Range("M62").Value = -7612.999899999999
# Row 65
$ws.Range("H65").Value = 2824.75
$ws.Range("I65").Value = 2766.3333
$ws.Range("K65").Value = 24896.9997
$ws.Range("M65").Value = -21464.9997
# Row 82
$ws.Range("H82").Value = 15715.214
# Row 85
$ws.Range("H85").Value = 15715.214
# Row 113
$ws.Range("H113").Value = 1229.7941
$ws.Range("J113").Value = 1249.963
$ws.Range("L113").Value = 3749.889
$ws.Range("N113").Value = -8089.889
# Row 137
$ws.Range("H137").Value = 8319.689
$ws.Range("J137").Value = 11007.235
$ws.Range("L137").Value = 33021.705
$ws.Range("N137").Value = -43221.705

$ws = $wb.Worksheets.Item("GSM")
# Row 46
$ws.Range("H46").Value = 20000
$ws.Range("I46").Value = 20000
$ws.Range("K46").Value = 20000
$ws.Range("M46").Value = -19844
# Row 102
$ws.Range("H102").Value = 4097.5
$ws.Range("I102").Value = 4037
$ws.Range("J102").Value = 4400
$ws.Range("K102").Value = 4037
$ws.Range("L102").Value = 4400
$ws.Range("M102").Value = -2415
$ws.Range("N102").Value = -7644
# Row 126
$ws.Range("H126").Value = 26306.666
$ws.Range("I126").Value = 36960
$ws.Range("K126").Value = 110880
$ws.Range("M126").Value = -108410
# Row 132
$ws.Range("H132").Value = 9289.553
$ws.Range("I132").Value = 7236.909
$ws.Range("K132").Value = 21710.727
$ws.Range("M132").Value = -19180.727

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3596.8
$ws.Range("I7").Value = 3596.8
$ws.Range("K7").Value = 3596.8
$ws.Range("M7").Value = -3484.8
# Row 46
$ws.Range("H46").Value = 3589.2104
$ws.Range("I46").Value = 733.3333
$ws.Range("K46").Value = 733.3333
$ws.Range("M46").Value = -545.3333
# Row 61
$ws.Range("H61").Value = 13782.643
$ws.Range("I61").Value = 14765.923
$ws.Range("K61").Value = 14765.923
$ws.Range("M61").Value = -14563.923
# Row 68
$ws.Range("H68").Value = 1999.8572
$ws.Range("I68").Value = 1999.8572
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1999.8572
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1250.8572
$ws.Range("N68").ClearContents()
# Row 71
$ws.Range("H71").Value = 1999.8572
$ws.Range("I71").Value = 1999.8572
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 9999.286
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -6255.286
$ws.Range("N71").ClearContents()
# Row 82
$ws.Range("H82").Value = 926.3103599999999
$ws.Range("I82").Value = 945.1070999999999
$ws.Range("J82").Value = 400
$ws.Range("K82").Value = 945.1070999999999
$ws.Range("L82").Value = 400
$ws.Range("M82").Value = -584.1070999999999
$ws.Range("N82").Value = -1122
# Row 85
$ws.Range("H85").Value = 926.3103599999999
$ws.Range("I85").Value = 945.1070999999999
$ws.Range("J85").Value = 400
$ws.Range("K85").Value = 945.1070999999999
$ws.Range("L85").Value = 400
$ws.Range("M85").Value = 302.8929000000001
$ws.Range("N85").Value = -2896
# Row 113
$ws.Range("H113").Value = 13782.643
$ws.Range("I113").Value = 14765.923
$ws.Range("K113").Value = 14765.923
$ws.Range("M113").Value = -12595.923
# Row 126
$ws.Range("H126").Value = 3596.8
$ws.Range("I126").Value = 3596.8
$ws.Range("K126").Value = 10790.4
$ws.Range("M126").Value = -8320.400000000001
# Row 132
$ws.Range("H132").Value = 6876695
$ws.Range("I132").Value = 9740417
$ws.Range("K132").Value = 29221251
$ws.Range("M132").Value = -29218721
